$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 408
$ws1.Range("F6").Value = 75
$ws1.Range("F8").Value = 61
$ws1.Range("F9").Value = 6722
$ws1.Range("F16").Value = 16073
$ws1.Range("F19").Value = 324
$ws1.Range("F23").Value = 5
$ws1.Range("F24").Value = 883

# Sheet "全部类型" (sheet4) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 408
$ws4.Range("F6").Value = 75
$ws4.Range("F9").Value = 61
$ws4.Range("F10").Value = 6722
$ws4.Range("F18").Value = 16073
$ws4.Range("F21").Value = 324
$ws4.Range("F27").Value = 5
$ws4.Range("F28").Value = 883
